$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "nowcast" (h=0) column of values to insert into column B (rows 2..16),
# shifting the existing Q0..Q9 data one column to the right. Any value that
# was already in column K falls off the fixed A:K grid when a row was full,
# matching the behavior of an insert-cells-shift-right operation.
$newValues = @{ 2 = 0.0000003593882045849206; 3 = -0.0000001035781544145298; 4 = -0.000000000387512216759589; 5 = -0.000000107388789361007; 6 = -0.0000001035472805832605; 7 = 0.000006303355340908645; 8 = -0.0000002375649628613696; 9 = 0.0000003720025918141356; 10 = 0.0000003829984367986761; 11 = -0.000003160475492397508; 12 = -0.00000004101096154340844; 13 = -0.0000001831659499074156; 14 = 0.0000002770877186031306; 15 = 0.000000229775004800814; 16 = -0.0000001554241066958895 }

for ($row = 2; $row -le 16; $row++) {
    # Shift existing values in this row one column to the right, within B:K,
    # dropping whatever would fall past column K.
    for ($col = 11; $col -ge 3; $col--) {
        $srcCell = $ws.Cells.Item($row, $col - 1)
        $dstCell = $ws.Cells.Item($row, $col)
        $dstCell.Value = $srcCell.Value2
    }
    # Place the new value into column B.
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
